$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.573.62"
$ws.Range("E2").Value = "  -4.42%  "
$ws.Range("D3").Value = "3.337.58"
$ws.Range("E3").Value = "  -1.36%  "
$ws.Range("D5").Value = "'573.39"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -3.34%  "
$ws.Range("D6").Value = "'180.73"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -5.23%  "
$ws.Range("D7").Value = "'0.625"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +2.62%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D10").Value = "'6.64"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.68%  "
$ws.Range("E11").Value = "  -3.63%  "
$ws.Range("D12").Value = "3.915.76"
$ws.Range("E12").Value = "  -1.47%  "
$ws.Range("E13").Value = "  -0.58%  "
$ws.Range("D14").Value = "'27.00"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -5.90%  "
$ws.Range("D15").Value = "66.705.60"
$ws.Range("E15").Value = "  -4.14%  "
$ws.Range("E16").Value = "  -2.47%  "
$ws.Range("D17").Value = "3.345.43"
$ws.Range("E17").Value = "  -0.92%  "
$ws.Range("D18").Value = "'436.31"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -3.34%  "
$ws.Range("D19").Value = "'5.68"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.39%  "
$ws.Range("D20").Value = "'13.55"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.83%  "
$ws.Range("D21").Value = "'7.59"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.77%  "
$ws.Range("D22").Value = "'73.54"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -3.36%  "
$ws.Range("D23").Value = "'1.00"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.20%  "
$ws.Range("D24").Value = "'0.516"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.89%  "
$ws.Range("D25").Value = "'0.0000117"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -3.85%  "
$ws.Range("D26").Value = "'0.190"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.50%  "
$ws.Range("E27").Value = "  -5.22%  "
$ws.Range("D28").Value = "'0.999"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.06%  "
$ws.Range("E29").Value = "  -3.23%  "
$ws.Range("D30").Value = "'22.81"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.76%  "
$ws.Range("B31").Value = "USDe"
$ws.Range("C31").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D31").Value = "'0.999"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.02%  "
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").Value = "'5.28"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -6.05%  "
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").Value = "'1.23"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -3.88%  "
$ws.Range("B34").Value = "Aptos"
$ws.Range("C34").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D34").Value = "'6.76"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -3.06%  "
$ws.Range("D35").Value = "'162.67"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.37%  "
$ws.Range("E36").Value = "  -5.86%  "
$ws.Range("D37").Value = "'27.50"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.07%  "
$ws.Range("D38").Value = "'1.80"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -7.50%  "
$ws.Range("D39").Value = "2.829.09"
$ws.Range("E39").Value = "  +3.02%  "
$ws.Range("D40").Value = "'0.797"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.16%  "
$ws.Range("E41").Value = "  -4.03%  "
$ws.Range("D42").Value = "'6.21"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -5.82%  "
$ws.Range("D43").Value = "'40.14"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.36%  "
$ws.Range("E44").Value = "  -3.12%  "
$ws.Range("D45").Value = "'24.48"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -3.81%  "
$ws.Range("E46").Value = "  -6.04%  "
$ws.Range("D47").Value = "'323.11"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -4.97%  "
$ws.Range("D48").Value = "'0.0273"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -3.90%  "
$ws.Range("D49").Value = "'0.102"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.65%  "
$ws.Range("D50").Value = "'0.980"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -3.79%  "
$ws.Range("E51").Value = "  -2.58%  "
